# University workbook fix-up:
# The backend that writes to this sheet ran twice for the same request,
# so duplicate rows / stale "second run" values were left behind. This
# script reverts the sheets back to the single-run state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# students
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("students")
$ws.Range("D2").Value = "hj"
$ws.Range("E2").Value = "h"
$ws.Range("F2").Value = "jh"
$ws.Range("G2").Value = "j"
$ws.Range("H2").Value = "hj"
$ws.Range("D3").Value = "hjh"
$ws.Range("E3").Value = "jhj"
$ws.Range("F3").Value = "h"
$ws.Range("G3").Value = "jh"
$ws.Range("J3").Value = 2
$ws.Range("A4:H4").ClearContents()
$ws.Activate()
$ws.Range("J3").Select()

# ---------------------------------------------------------------
# student_courses
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("student_courses")
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Hasith"
$ws.Range("C3").Value = "Dewmina"
$ws.Range("A4:E4").ClearContents()
$ws.Range("L4").Value = 2
$ws.Range("A5:E6").ClearContents()

# ---------------------------------------------------------------
# student_pswd
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("student_pswd")
$ws.Range("A4:C4").ClearContents()
$ws.Range("G6").Value = 2

# ---------------------------------------------------------------
# Teachers
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Teachers")
$ws.Range("B2").Value = "Sam"
$ws.Range("C2").Value = "j"
$ws.Range("D2").Value = "kjk"
$ws.Range("E2").Value = "j"
$ws.Range("A3:E3").ClearContents()
$ws.Range("H4").Value = 1
$ws.Activate()
$ws.Range("H4").Select()

# ---------------------------------------------------------------
# teacher_psswd
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("teacher_psswd")
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "lecturer1"
$ws.Range("C3").Value = "Sam"

# ---------------------------------------------------------------
# notifications
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("notifications")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "COURSE ENROLLMENT"
$ws.Range("C2").Value = "Hasith"
$ws.Range("D2").Value = "Matt"
$ws.Range("E2").Value = "Hasith would like to enroll in the Maths"
$ws.Range("A3:E3").ClearContents()
$ws.Range("A4:H4").ClearContents()

# Restore original active sheet/tab (notifications, unchanged by this fix)
$ws.Activate()
$ws.Range("H3").Select()
